# Regenerate save_data to use K instead of Strike#, recalculated std/mean,
# recalculated and written s_vals. This updates column G ("K") values for
# the affected rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 3
    13 = 1
    15 = 2
    16 = 2
    17 = 3
    18 = 0
    19 = 2
    20 = 2
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
